$d = $word.ActiveDocument

# --- Locate the paragraph ending in "...nie zostanie wykonana żadna zamiana."
#     and insert a brand-new empty paragraph right after it ------------------
$anchor = $d.Content
$anchor.Find.Execute("żadna zamiana.", $false, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()

# --- Re-resolve the paragraph collection (the handle used above can go stale
#     once the document structure changes) and find the new, empty paragraph
#     that now immediately follows the anchor paragraph. ---------------------
$paragraphs = $d.Paragraphs
$count = $paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    if ($paragraphs.Item($i).Range.Text.Contains("wykonana")) {
        $anchorIndex = $i
        break
    }
}

$newParaRange = $paragraphs.Item($anchorIndex + 1).Range

# --- Fill the new paragraph with the "Złożoność obliczeniowa O(n^2)." text,
#     expressed as an inline OOXML math equation. -----------------------------
$mathXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:r><w:t xml:space="preserve">Złożoność obliczeniowa </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>O(</m:t></m:r><m:sSup><m:sSupPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSupPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>n</m:t></m:r></m:e><m:sup><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>)</m:t></m:r></m:oMath><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>.</w:t></w:r></w:p>'

$newParaRange.InsertXML($mathXml)
